# Auto-generated edit script: apply numeric value changes described by the diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4353.5713
$ws.Range("I10").Value = 1983.3334
$ws.Range("K10").Value = 1983.3334
$ws.Range("M10").Value = -1690.3334
$ws.Range("H21").Value = 16678.438
$ws.Range("I21").Value = 17046.25
$ws.Range("K21").Value = 17046.25
$ws.Range("M21").Value = -16578.25
$ws.Range("H23").Value = 16678.438
$ws.Range("I23").Value = 17046.25
$ws.Range("K23").Value = 17046.25
$ws.Range("M23").Value = -16812.25
$ws.Range("H62").Value = 2196.4614
$ws.Range("I62").Value = 2328
$ws.Range("J62").Value = 2017.091
$ws.Range("K62").Value = 2328
$ws.Range("L62").Value = 2017.091
$ws.Range("M62").Value = -1704
$ws.Range("N62").Value = -3265.091
$ws.Range("H64").Value = 3718
$ws.Range("I64").Value = 3718
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3718
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3470
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 2196.4614
$ws.Range("I65").Value = 2328
$ws.Range("J65").Value = 2017.091
$ws.Range("K65").Value = 11640
$ws.Range("L65").Value = 10085.455
$ws.Range("M65").Value = -8520
$ws.Range("N65").Value = -16325.455
$ws.Range("H67").Value = 3718
$ws.Range("I67").Value = 3718
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3718
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2860
$ws.Range("N67").ClearContents()
$ws.Range("H138").Value = 4935.7065
$ws.Range("J138").Value = 5268.4814
$ws.Range("L138").Value = 15805.4442
$ws.Range("N138").Value = -26085.4442

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 20261.82
$ws.Range("I32").Value = 15893.299
$ws.Range("J32").Value = 78800
$ws.Range("K32").Value = 15893.299
$ws.Range("L32").Value = 78800
$ws.Range("M32").Value = -15606.299
$ws.Range("N32").Value = -79374
$ws.Range("H37").Value = 10123.777
$ws.Range("J37").Value = 10123.777
$ws.Range("L37").Value = 10123.777
$ws.Range("N37").Value = -10669.777
$ws.Range("H74").Value = 9094668
$ws.Range("I74").Value = 13889635
$ws.Range("J74").Value = 9466.315000000001
$ws.Range("K74").Value = 13889635
$ws.Range("L74").Value = 9466.315000000001
$ws.Range("M74").Value = -13888761
$ws.Range("N74").Value = -11214.315
$ws.Range("H77").Value = 9094668
$ws.Range("I77").Value = 13889635
$ws.Range("J77").Value = 9466.315000000001
$ws.Range("K77").Value = 69448175
$ws.Range("L77").Value = 47331.575
$ws.Range("M77").Value = -69443807
$ws.Range("N77").Value = -56067.575

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 34399.668
$ws.Range("J6").Value = 34399.668
$ws.Range("L6").Value = 34399.668
$ws.Range("N6").Value = -34625.668
$ws.Range("H27").Value = 39542
$ws.Range("J27").Value = 39542
$ws.Range("L27").Value = 39542
$ws.Range("N27").Value = -39926
$ws.Range("H82").Value = 10030.7
$ws.Range("I82").Value = 2884.5
$ws.Range("J82").Value = 20750
$ws.Range("K82").Value = 2884.5
$ws.Range("L82").Value = 20750
$ws.Range("M82").Value = -2501.5
$ws.Range("N82").Value = -21516
$ws.Range("H85").Value = 10030.7
$ws.Range("I85").Value = 2884.5
$ws.Range("J85").Value = 20750
$ws.Range("K85").Value = 2884.5
$ws.Range("L85").Value = 20750
$ws.Range("M85").Value = -1558.5
$ws.Range("N85").Value = -23402
$ws.Range("H99").Value = 1299.1
$ws.Range("I99").Value = 1223.4445
$ws.Range("J99").Value = 1980
$ws.Range("K99").Value = 1223.4445
$ws.Range("L99").Value = 1980
$ws.Range("M99").Value = 274.5554999999999
$ws.Range("N99").Value = -4976
$ws.Range("H134").Value = 8697731
$ws.Range("I134").Value = 8697731
$ws.Range("K134").Value = 26093193
$ws.Range("M134").Value = -26090658

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1958.1818
$ws.Range("I19").Value = 171.11111
$ws.Range("K19").Value = 171.11111
$ws.Range("M19").Value = -1.111109999999996
$ws.Range("H24").Value = 1958.1818
$ws.Range("I24").Value = 171.11111
$ws.Range("K24").Value = 171.11111
$ws.Range("M24").Value = -1.111109999999996
$ws.Range("H31").Value = 13948.151
$ws.Range("I31").Value = 19698.408
$ws.Range("J31").Value = 7976.731
$ws.Range("K31").Value = 19698.408
$ws.Range("L31").Value = 7976.731
$ws.Range("M31").Value = -19403.408
$ws.Range("N31").Value = -8566.731
$ws.Range("H34").Value = 13948.151
$ws.Range("I34").Value = 19698.408
$ws.Range("J34").Value = 7976.731
$ws.Range("K34").Value = 19698.408
$ws.Range("L34").Value = 7976.731
$ws.Range("M34").Value = -19496.408
$ws.Range("N34").Value = -8380.731
$ws.Range("H86").Value = 2113.04
$ws.Range("I86").Value = 2233.6
$ws.Range("J86").Value = 1932.2
$ws.Range("K86").Value = 2233.6
$ws.Range("L86").Value = 1932.2
$ws.Range("M86").Value = -1110.6
$ws.Range("N86").Value = -4178.2
$ws.Range("H89").Value = 2113.04
$ws.Range("I89").Value = 2233.6
$ws.Range("J89").Value = 1932.2
$ws.Range("K89").Value = 11168
$ws.Range("L89").Value = 9661
$ws.Range("M89").Value = -5552
$ws.Range("N89").Value = -20893
$ws.Range("H122").Value = 12690.223
$ws.Range("I122").Value = 14089
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 42267
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -39817
$ws.Range("N122").Value = -9400

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 320.5
$ws.Range("I21").Value = 320.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 961.5
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -788.5
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 120.26316
$ws.Range("I23").Value = 46.666668
$ws.Range("J23").Value = 134.0625
$ws.Range("K23").Value = 140.000004
$ws.Range("L23").Value = 402.1875
$ws.Range("M23").Value = 94.99999600000001
$ws.Range("N23").Value = -872.1875
$ws.Range("H107").Value = 125308.69
$ws.Range("I107").Value = 91208.63
$ws.Range("J107").Value = 200328.8
$ws.Range("K107").Value = 273625.89
$ws.Range("L107").Value = 600986.3999999999
$ws.Range("M107").Value = -271705.89
$ws.Range("N107").Value = -604826.3999999999
$ws.Range("H109").Value = 2599.7693
$ws.Range("I109").Value = 1399.5
$ws.Range("K109").Value = 4198.5
$ws.Range("M109").Value = -3158.5
$ws.Range("H115").Value = 6023.3335
$ws.Range("I115").Value = 1340
$ws.Range("K115").Value = 4020
$ws.Range("M115").Value = -2845
$ws.Range("H122").Value = 5953528
$ws.Range("I122").Value = 10753210
$ws.Range("J122").Value = 822833.4399999999
$ws.Range("K122").Value = 96778890
$ws.Range("L122").Value = 7405500.959999999
$ws.Range("M122").Value = -96776440
$ws.Range("N122").Value = -7410400.959999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1146.2162
$ws.Range("I122").Value = 1038.2333
$ws.Range("J122").Value = 1609
$ws.Range("K122").Value = 3114.699900000001
$ws.Range("L122").Value = 4827
$ws.Range("M122").Value = -664.6999000000005
$ws.Range("N122").Value = -9727
$ws.Range("H132").Value = 3039.389
$ws.Range("I132").Value = 3802.158
$ws.Range("J132").Value = 2186.8823
$ws.Range("K132").Value = 11406.474
$ws.Range("L132").Value = 6560.646900000001
$ws.Range("M132").Value = -8876.474
$ws.Range("N132").Value = -11620.6469

